# Insert a new data row at row 824 (pushing the existing rows 824..863 down
# to 825..864), then fill the newly inserted row with the new record's
# values. Excel's row Insert() naturally carries the formatting (e.g. the
# date style on column D) down from the row immediately above, which matches
# the original author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 824-863 down to 825-864 by inserting a fresh row at 824.
$ws.Rows.Item(824).Insert()

# Populate the newly inserted row 824 with the new record.
$ws.Cells.Item(824, 1).Value  = 3
$ws.Cells.Item(824, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(824, 3).Value  = "Coquimbo"
$ws.Cells.Item(824, 4).Value  = 45267
$ws.Cells.Item(824, 5).Value  = 5
$ws.Cells.Item(824, 6).Value  = 100112037
$ws.Cells.Item(824, 7).Value  = "Cebollín"
$ws.Cells.Item(824, 8).Value  = "Sin especificar"
$ws.Cells.Item(824, 9).Value  = "Primera"
$ws.Cells.Item(824, 10).Value = 220
$ws.Cells.Item(824, 11).Value = 3500
$ws.Cells.Item(824, 12).Value = 4000
$ws.Cells.Item(824, 13).Value = 3773
$ws.Cells.Item(824, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(824, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(824, 16).Value = 105
$ws.Cells.Item(824, 17).Value = 36
$ws.Cells.Item(824, 18).Value = "Hortaliza"
